$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 734:735. Excel shifts the existing rows 734+ down by 2,
# so the old row 734 data (previously date 44616 / Pintón) ends up at row 736,
# and the old row 735 data (previously date 44616 / Primera Pintón) ends up at
# row 737. Everything below keeps moving down by 2 all the way to the old
# last row 840, which becomes 842.
$ws.Rows("734:735").Insert()

# Populate the newly-inserted row 734 ("Pintón" bucket) with the new
# observation (same dimensional attributes as the rest of the sheet).
$ws.Range("A734").Value = 7
$ws.Range("B734").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C734").Value = "Ñuble"
$ws.Range("D734").Value = 44984
$ws.Range("E734").Value = 16
$ws.Range("F734").Value = "Fruta"
$ws.Range("G734").Value = 100108
$ws.Range("H734").Value = "Tropicales y subtropicales"
$ws.Range("I734").Value = 100108006
$ws.Range("J734").Value = "Plátano"
$ws.Range("K734").Value = "Sin especificar"
$ws.Range("L734").Value = "Pintón"
$ws.Range("M734").Value = 300
$ws.Range("N734").Value = 23000
$ws.Range("O734").Value = 23000
$ws.Range("P734").Value = 23000
$ws.Range("Q734").Value = "$/caja 20 kilos"
$ws.Range("R734").Value = "Ecuador"
$ws.Range("S734").Value = 1150
$ws.Range("T734").Value = 20

# Populate the newly-inserted row 735 ("Primera Pintón" bucket).
$ws.Range("A735").Value = 7
$ws.Range("B735").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C735").Value = "Ñuble"
$ws.Range("D735").Value = 44984
$ws.Range("E735").Value = 16
$ws.Range("F735").Value = "Fruta"
$ws.Range("G735").Value = 100108
$ws.Range("H735").Value = "Tropicales y subtropicales"
$ws.Range("I735").Value = 100108006
$ws.Range("J735").Value = "Plátano"
$ws.Range("K735").Value = "Sin especificar"
$ws.Range("L735").Value = "Primera Pintón"
$ws.Range("M735").Value = 500
$ws.Range("N735").Value = 24000
$ws.Range("O735").Value = 25000
$ws.Range("P735").Value = 24500
$ws.Range("Q735").Value = "$/caja 20 kilos"
$ws.Range("R735").Value = "Ecuador"
$ws.Range("S735").Value = 1225
$ws.Range("T735").Value = 20

# Make sure the date cells keep the date/time number format used throughout
# column D (style index 2 in the original file -> numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D734").NumberFormat = $ws.Range("D736").NumberFormat
$ws.Range("D735").NumberFormat = $ws.Range("D736").NumberFormat
